$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.439.06'
$ws.Range('E2').Value = '  -6.06%  '
$ws.Range('D3').Value = '2.527.34'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('D4').Formula = "'0.999"
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Formula = "'297.23"
$ws.Range('E5').Value = '  -3.36%  '
$ws.Range('D6').Formula = "'94.86"
$ws.Range('E6').Value = '  -5.33%  '
$ws.Range('D7').Formula = "'0.577"
$ws.Range('E7').Value = '  -4.28%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Formula = "'0.553"
$ws.Range('E9').Value = '  -4.58%  '
$ws.Range('D10').Formula = "'36.72"
$ws.Range('E10').Value = '  -6.98%  '
$ws.Range('D11').Formula = "'0.0807"
$ws.Range('E11').Value = '  -4.36%  '
$ws.Range('D12').Formula = "'7.72"
$ws.Range('E12').Value = '  -5.32%  '
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '2.909.76'
$ws.Range('E14').Value = '  -3.23%  '
$ws.Range('D15').Value = '2.524.77'
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('D16').Formula = "'0.874"
$ws.Range('E16').Value = '  -4.85%  '
$ws.Range('D17').Formula = "'14.15"
$ws.Range('E17').Value = '  -5.13%  '
$ws.Range('D18').Value = '43.433.66'
$ws.Range('E18').Value = '  -6.45%  '
$ws.Range('D19').Value = '0.0₃0970'
$ws.Range('E19').Value = '  -4.14%  '
$ws.Range('D20').Formula = "'6.57"
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('D21').Formula = "'12.43"
$ws.Range('E21').Value = '  -3.76%  '
$ws.Range('D22').Formula = "'72.39"
$ws.Range('E22').Value = '  +1.44%  '
$ws.Range('D23').Formula = "'261.29"
$ws.Range('E23').Value = '  -4.22%  '
$ws.Range('E24').Value = '  -4.80%  '
$ws.Range('D25').Formula = "'2.15"
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('D26').Formula = "'28.91"
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('D27').Formula = "'1.00"
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').Formula = "'10.07"
$ws.Range('E28').Value = '  -4.73%  '
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('D30').Formula = "'37.26"
$ws.Range('E30').Value = '  -4.21%  '
$ws.Range('D31').Formula = "'6.09"
$ws.Range('E31').Value = '  -3.87%  '
$ws.Range('D32').Formula = "'3.47"
$ws.Range('E32').Value = '  -4.98%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Formula = "'150.58"
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Formula = "'2.75"
$ws.Range('E34').Value = '  -3.85%  '
$ws.Range('D35').Formula = "'2.15"
$ws.Range('E35').Value = '  -3.53%  '
$ws.Range('D36').Formula = "'0.0802"
$ws.Range('E36').Value = '  -4.40%  '
$ws.Range('E37').Value = '  -5.71%  '
$ws.Range('D38').Formula = "'0.120"
$ws.Range('E38').Value = '  -3.03%  '
$ws.Range('D39').Formula = "'23.44"
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('D40').Formula = "'16.23"
$ws.Range('E40').Value = '  +2.06%  '
$ws.Range('D41').Formula = "'3.51"
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('D42').Formula = "'0.0310"
$ws.Range('E42').Value = '  -6.19%  '
$ws.Range('D43').Formula = "'3.82"
$ws.Range('E43').Value = '  -5.81%  '
$ws.Range('D44').Value = '2.018.36'
$ws.Range('E44').Value = '  -4.56%  '
$ws.Range('D45').Formula = "'0.997"
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Formula = "'86.04"
$ws.Range('E46').Value = '  -7.57%  '
$ws.Range('D47').Formula = "'1.61"
$ws.Range('E47').Value = '  +3.14%  '
$ws.Range('D48').Formula = "'8.93"
$ws.Range('E48').Value = '  -6.07%  '
$ws.Range('D49').Value = '2.768.34'
$ws.Range('E49').Value = '  -3.32%  '
$ws.Range('D50').Formula = "'103.57"
$ws.Range('E50').Value = '  -4.90%  '
$ws.Range('D51').Formula = "'0.189"
$ws.Range('E51').Value = '  -5.76%  '
